$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("A4").Value = 131257424
$ws.Range("B4").Value = 79245
$ws.Range("E4").Value = 6425
$ws.Range("F4").Value = "Garnlav"
$ws.Range("G4").Value = "Alectoria sarmentosa"
$ws.Range("H4").Value = "(Ach.) Ach."
$ws.Range("M4").ClearContents()
$ws.Range("Q4").Value = 488876
$ws.Range("R4").Value = 6665177
$ws.Range("Z4").Value = "11:33"
$ws.Range("AB4").Value = "11:33"
$ws.Range("AC4").Value = "Gran"

# Row 5
$ws.Range("A5").Value = 131255793
$ws.Range("B5").Value = 91830
$ws.Range("E5").Value = 5432
$ws.Range("F5").Value = "Granticka"
$ws.Range("G5").Value = "Porodaedalea chrysoloma s.lat."
$ws.Range("H5").ClearContents()
$ws.Range("M5").ClearContents()
$ws.Range("Q5").Value = 488817
$ws.Range("R5").Value = 6665110
$ws.Range("Z5").Value = "09:56"
$ws.Range("AB5").Value = "09:56"
$ws.Range("AC5").Value = "Flera fruktkroppar."

# Row 6
$ws.Range("A6").Value = 131256691
$ws.Range("B6").Value = 57884
$ws.Range("E6").Value = 100109
$ws.Range("F6").Value = "Tretåig hackspett"
$ws.Range("G6").Value = "Picoides tridactylus"
$ws.Range("H6").Value = "(Linnaeus, 1758)"
$ws.Range("M6").Value = "äldre spår"
$ws.Range("Q6").Value = 488667
$ws.Range("R6").Value = 6665262
$ws.Range("Z6").Value = "10:55"
$ws.Range("AB6").Value = "10:55"
$ws.Range("AC6").Value = "Ringhack på gran."

# Row 7
$ws.Range("A7").Value = 131260583
$ws.Range("B7").Value = 57884
$ws.Range("E7").Value = 100109
$ws.Range("F7").Value = "Tretåig hackspett"
$ws.Range("G7").Value = "Picoides tridactylus"
$ws.Range("H7").Value = "(Linnaeus, 1758)"
$ws.Range("M7").Value = "färska spår"
$ws.Range("Q7").Value = 488834
$ws.Range("R7").Value = 6665228
$ws.Range("Z7").Value = "15:30"
$ws.Range("AB7").Value = "15:30"
$ws.Range("AC7").Value = "Ringhack på tall."

# Row 11
$ws.Range("A11").Value = 131257290
$ws.Range("B11").Value = 57884
$ws.Range("E11").Value = 100109
$ws.Range("F11").Value = "Tretåig hackspett"
$ws.Range("G11").Value = "Picoides tridactylus"
$ws.Range("H11").Value = "(Linnaeus, 1758)"
$ws.Range("M11").Value = "äldre spår"
$ws.Range("Q11").Value = 488842
$ws.Range("R11").Value = 6665224
$ws.Range("Z11").Value = "11:26"
$ws.Range("AB11").Value = "11:26"
$ws.Range("AC11").Value = "Ringhack på tall."

# Row 12
$ws.Range("A12").Value = 131256673
$ws.Range("Q12").Value = 488652
$ws.Range("R12").Value = 6665282
$ws.Range("Z12").Value = "10:54"
$ws.Range("AB12").Value = "10:54"
$ws.Range("AC12").Value = "Ringhack på tall."

# Row 13
$ws.Range("A13").Value = 131257520
$ws.Range("B13").Value = 79245
$ws.Range("E13").Value = 6425
$ws.Range("F13").Value = "Garnlav"
$ws.Range("G13").Value = "Alectoria sarmentosa"
$ws.Range("H13").Value = "(Ach.) Ach."
$ws.Range("M13").ClearContents()
$ws.Range("Q13").Value = 488939
$ws.Range("R13").Value = 6665149
$ws.Range("Z13").Value = "11:41"
$ws.Range("AB13").Value = "11:41"
$ws.Range("AC13").Value = "Gran"

# Row 14
$ws.Range("A14").Value = 131260641
$ws.Range("Q14").Value = 488859
$ws.Range("R14").Value = 6665292
$ws.Range("Z14").Value = "15:34"
$ws.Range("AB14").Value = "15:34"
$ws.Range("AC14").Value = "Ringhack på gran."

# Row 35
$ws.Range("A35").Value = 131260531
$ws.Range("B35").Value = 79245
$ws.Range("E35").Value = 6425
$ws.Range("F35").Value = "Garnlav"
$ws.Range("G35").Value = "Alectoria sarmentosa"
$ws.Range("H35").Value = "(Ach.) Ach."
$ws.Range("Q35").Value = 488786
$ws.Range("R35").Value = 6665188
$ws.Range("Z35").Value = "15:25"
$ws.Range("AB35").Value = "15:25"
$ws.Range("AC35").Value = "Gran"

# Row 36
$ws.Range("A36").Value = 131257385
$ws.Range("B36").Value = 91830
$ws.Range("E36").Value = 5432
$ws.Range("F36").Value = "Granticka"
$ws.Range("G36").Value = "Porodaedalea chrysoloma s.lat."
$ws.Range("H36").ClearContents()
$ws.Range("Q36").Value = 488876
$ws.Range("R36").Value = 6665194
$ws.Range("Z36").Value = "11:31"
$ws.Range("AB36").Value = "11:31"
$ws.Range("AC36").Value = "Lågstubbe."

# Row 42
$ws.Range("A42").Value = 131273946
$ws.Range("Q42").Value = 488774
$ws.Range("R42").Value = 6665353

# Row 43
$ws.Range("A43").Value = 131273991
$ws.Range("Q43").Value = 488928
$ws.Range("R43").Value = 6665146
